# Scheduled-runner data refresh: pulls fresh Universalis market prices
# and rewrites the derived price/profit columns (H:N) for the affected
# leve rows on each job sheet. Row/column layout is unchanged; only the
# computed numbers move.
#   H currentAveragePrice    K LevePriceNQ
#   I currentAveragePriceNQ  L LevePriceHQ
#   J currentAveragePriceHQ  M LeveProfitNQ
#                            N LeveProfitHQ
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 86: Filling in the Blanks
$ws.Cells.Item(86, 8).Value = 158529.77
$ws.Cells.Item(86, 9).Value = 128874.125
$ws.Cells.Item(86, 10).Value = 205978.8
$ws.Cells.Item(86, 11).Value = 128874.125
$ws.Cells.Item(86, 12).Value = 205978.8
$ws.Cells.Item(86, 13).Value = -127751.125
$ws.Cells.Item(86, 14).Value = -208224.8

# Row 89: Ink into Antiquity (L)
$ws.Cells.Item(89, 8).Value = 158529.77
$ws.Cells.Item(89, 9).Value = 128874.125
$ws.Cells.Item(89, 10).Value = 205978.8
$ws.Cells.Item(89, 11).Value = 644370.625
$ws.Cells.Item(89, 12).Value = 1029894
$ws.Cells.Item(89, 13).Value = -638754.625
$ws.Cells.Item(89, 14).Value = -1041126

# Row 134: Binding Spells
$ws.Cells.Item(134, 8).Value = 68535.71000000001
$ws.Cells.Item(134, 10).Value = 68535.71000000001
$ws.Cells.Item(134, 12).Value = 68535.71000000001
$ws.Cells.Item(134, 14).Value = -78675.71000000001

# Row 137: Cutting Edge of Culinary Quality
$ws.Cells.Item(137, 8).Value = 3194.509
$ws.Cells.Item(137, 9).Value = 1851.579
$ws.Cells.Item(137, 10).Value = 3903.2778
$ws.Cells.Item(137, 11).Value = 5554.737
$ws.Cells.Item(137, 12).Value = 11709.8334
$ws.Cells.Item(137, 13).Value = -3004.737
$ws.Cells.Item(137, 14).Value = -16809.8334

# Row 138: All-night Crafting
$ws.Cells.Item(138, 8).Value = 1471662.5
$ws.Cells.Item(138, 9).Value = 4450.4443
$ws.Cells.Item(138, 10).Value = 1720811.8
$ws.Cells.Item(138, 11).Value = 13351.3329
$ws.Cells.Item(138, 12).Value = 5162435.4
$ws.Cells.Item(138, 13).Value = -8211.332900000001
$ws.Cells.Item(138, 14).Value = -5172715.4

$ws = $wb.Worksheets.Item("ARM")
# Row 32: Ingot We Trust
$ws.Cells.Item(32, 8).Value = 15003.527
$ws.Cells.Item(32, 9).Value = 13128.969
$ws.Cells.Item(32, 11).Value = 13128.969
$ws.Cells.Item(32, 13).Value = -12841.969

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin
$ws.Cells.Item(86, 8).Value = 2515.5862
$ws.Cells.Item(86, 9).Value = 2649.5652
$ws.Cells.Item(86, 10).Value = 2002
$ws.Cells.Item(86, 11).Value = 2649.5652
$ws.Cells.Item(86, 12).Value = 2002
$ws.Cells.Item(86, 13).Value = -1526.5652
$ws.Cells.Item(86, 14).Value = -4248

# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Cells.Item(89, 8).Value = 2515.5862
$ws.Cells.Item(89, 9).Value = 2649.5652
$ws.Cells.Item(89, 10).Value = 2002
$ws.Cells.Item(89, 11).Value = 13247.826
$ws.Cells.Item(89, 12).Value = 10010
$ws.Cells.Item(89, 13).Value = -7631.826000000001
$ws.Cells.Item(89, 14).Value = -21242

# Row 134: Ruthenium Supremium
$ws.Cells.Item(134, 8).Value = 23830.709
$ws.Cells.Item(134, 9).Value = 2800.5293
$ws.Cells.Item(134, 11).Value = 8401.5879
$ws.Cells.Item(134, 13).Value = -5866.5879

$ws = $wb.Worksheets.Item("CUL")
# Row 80: Saucy for a Suitor
$ws.Cells.Item(80, 8).Value = 1368.3334
$ws.Cells.Item(80, 9).Value = 1155
$ws.Cells.Item(80, 10).Value = 1475
$ws.Cells.Item(80, 11).Value = 3465
$ws.Cells.Item(80, 12).Value = 4425
$ws.Cells.Item(80, 13).Value = -2529
$ws.Cells.Item(80, 14).Value = -6297

# Row 83: Saved by the Sauce (L)
$ws.Cells.Item(83, 8).Value = 1368.3334
$ws.Cells.Item(83, 9).Value = 1155
$ws.Cells.Item(83, 10).Value = 1475
$ws.Cells.Item(83, 11).Value = 10395
$ws.Cells.Item(83, 12).Value = 13275
$ws.Cells.Item(83, 13).Value = -5715
$ws.Cells.Item(83, 14).Value = -22635

# Row 122: Salt of the North
$ws.Cells.Item(122, 8).Value = 1303.2667
$ws.Cells.Item(122, 9).Value = 616.86664
$ws.Cells.Item(122, 10).Value = 1646.4667
$ws.Cells.Item(122, 11).Value = 5551.79976
$ws.Cells.Item(122, 12).Value = 14818.2003
$ws.Cells.Item(122, 13).Value = -3101.79976
$ws.Cells.Item(122, 14).Value = -19718.2003

# Row 131: The Mountain Steeped
$ws.Cells.Item(131, 8).Value = 511.98914
$ws.Cells.Item(131, 9).Value = 285.43103
$ws.Cells.Item(131, 10).Value = 898.4706
$ws.Cells.Item(131, 11).Value = 856.2930900000001
$ws.Cells.Item(131, 12).Value = 2695.4118
$ws.Cells.Item(131, 13).Value = 4183.70691
$ws.Cells.Item(131, 14).Value = -12775.4118

$ws = $wb.Worksheets.Item("GSM")
# Row 2: Copper and Robbers
$ws.Cells.Item(2, 8).Value = 328.125
$ws.Cells.Item(2, 10).Value = 555
$ws.Cells.Item(2, 12).Value = 555
$ws.Cells.Item(2, 14).Value = -781

# Row 18: Gorgeous Gorget
$ws.Cells.Item(18, 8).Value = 10000
$ws.Cells.Item(18, 10).Value = 10000
$ws.Cells.Item(18, 12).Value = 10000
$ws.Cells.Item(18, 14).Value = -10586

# Row 43: Get the Green Stuff
$ws.Cells.Item(43, 8).Value = 3989
$ws.Cells.Item(43, 9).Value = 1786.8
$ws.Cells.Item(43, 11).Value = 1786.8
$ws.Cells.Item(43, 13).Value = -1635.8

# Row 46: Burning the Midnight Oil
$ws.Cells.Item(46, 8).Value = 16360.667

# Row 57: Gold Is So Last Year
$ws.Cells.Item(57, 8).Value = 14000
$ws.Cells.Item(57, 9).Value = 2000
$ws.Cells.Item(57, 10).Value = 20000
$ws.Cells.Item(57, 11).Value = 2000
$ws.Cells.Item(57, 12).Value = 20000
$ws.Cells.Item(57, 13).Value = -1180
$ws.Cells.Item(57, 14).Value = -21640

# Row 80: Needs More Prayerbell
$ws.Cells.Item(80, 8).Value = 6043.2173
$ws.Cells.Item(80, 9).Value = 11275
$ws.Cells.Item(80, 10).Value = 4196.706
$ws.Cells.Item(80, 11).Value = 11275
$ws.Cells.Item(80, 12).Value = 4196.706
$ws.Cells.Item(80, 13).Value = -10277
$ws.Cells.Item(80, 14).Value = -6192.706

# Row 83: With a Noise That Reaches Heaven (L)
$ws.Cells.Item(83, 8).Value = 6043.2173
$ws.Cells.Item(83, 9).Value = 11275
$ws.Cells.Item(83, 10).Value = 4196.706
$ws.Cells.Item(83, 11).Value = 56375
$ws.Cells.Item(83, 12).Value = 20983.53
$ws.Cells.Item(83, 13).Value = -51383
$ws.Cells.Item(83, 14).Value = -30967.53

# Row 93: One Ring Circus
$ws.Cells.Item(93, 8).Value = 28975.5
$ws.Cells.Item(93, 10).Value = 28975.5
$ws.Cells.Item(93, 12).Value = 28975.5
$ws.Cells.Item(93, 14).Value = -32719.5

$ws = $wb.Worksheets.Item("LTW")
# Row 136: Respect for Br'aax
$ws.Cells.Item(136, 8).Value = 4427.4683
$ws.Cells.Item(136, 9).Value = 2702
$ws.Cells.Item(136, 10).Value = 8494.643
$ws.Cells.Item(136, 11).Value = 8106
$ws.Cells.Item(136, 12).Value = 25483.929
$ws.Cells.Item(136, 13).Value = -5556
$ws.Cells.Item(136, 14).Value = -30583.929

$ws = $wb.Worksheets.Item("WVR")
# Row 45: Private Concerns
$ws.Cells.Item(45, 8).Value = 6276.5
$ws.Cells.Item(45, 9).Value = 2480
$ws.Cells.Item(45, 10).Value = 7542
$ws.Cells.Item(45, 11).Value = 2480
$ws.Cells.Item(45, 12).Value = 7542
$ws.Cells.Item(45, 13).Value = -1989
$ws.Cells.Item(45, 14).Value = -8524

# Row 47: The Wages of Sin
$ws.Cells.Item(47, 8).Value = 174500
$ws.Cells.Item(47, 10).Value = 174500
$ws.Cells.Item(47, 12).Value = 174500
$ws.Cells.Item(47, 14).Value = -175644

# Row 81: Where the Dragonflies, the Net Catches
$ws.Cells.Item(81, 8).Value = 4300.8
$ws.Cells.Item(81, 9).Value = 3500
$ws.Cells.Item(81, 10).Value = 4501
$ws.Cells.Item(81, 11).Value = 7000
$ws.Cells.Item(81, 12).Value = 9002
$ws.Cells.Item(81, 13).Value = -5939
$ws.Cells.Item(81, 14).Value = -11124

# Row 84: To Kill a Dragon on Nameday (L)
$ws.Cells.Item(84, 8).Value = 4300.8
$ws.Cells.Item(84, 9).Value = 3500
$ws.Cells.Item(84, 10).Value = 4501
$ws.Cells.Item(84, 11).Value = 35000
$ws.Cells.Item(84, 12).Value = 45010
$ws.Cells.Item(84, 13).Value = -29696
$ws.Cells.Item(84, 14).Value = -55618

# Row 122: Heavy Armoire
$ws.Cells.Item(122, 8).Value = 2892.2415
$ws.Cells.Item(122, 9).Value = 1114.75
$ws.Cells.Item(122, 10).Value = 5079.923
$ws.Cells.Item(122, 11).Value = 3344.25
$ws.Cells.Item(122, 12).Value = 15239.769
$ws.Cells.Item(122, 13).Value = -894.25
$ws.Cells.Item(122, 14).Value = -20139.769

# Row 136: Weaving the Envelope
$ws.Cells.Item(136, 8).Value = 4534.6665
$ws.Cells.Item(136, 9).Value = 2382.4783
$ws.Cells.Item(136, 10).Value = 7628.4375
$ws.Cells.Item(136, 11).Value = 7147.4349
$ws.Cells.Item(136, 12).Value = 22885.3125
$ws.Cells.Item(136, 13).Value = -4597.4349
$ws.Cells.Item(136, 14).Value = -27985.3125
